$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 768.625
$ws.Range("I12").Value = 724
$ws.Range("J12").Value = 775
$ws.Range("K12").Value = 724
$ws.Range("L12").Value = 775
$ws.Range("M12").Value = -554
$ws.Range("N12").Value = -1115

$ws.Range("H19").Value = 295.91666
$ws.Range("I19").Value = 273.33334
$ws.Range("J19").Value = 363.66666
$ws.Range("K19").Value = 273.33334
$ws.Range("L19").Value = 363.66666
$ws.Range("M19").Value = -98.33334000000002
$ws.Range("N19").Value = -713.66666

$ws.Range("H33").Value = 730.125
$ws.Range("I33").Value = 731.3333
$ws.Range("J33").Value = 729.4
$ws.Range("K33").Value = 731.3333
$ws.Range("L33").Value = 729.4
$ws.Range("M33").Value = -502.3333
$ws.Range("N33").Value = -1187.4

$ws.Range("H43").Value = 2240.6155
$ws.Range("I43").Value = 2572.75
$ws.Range("K43").Value = 2572.75
$ws.Range("M43").Value = -2503.75

$ws.Range("H58").Value = 1800
$ws.Range("I58").Value = 1800
$ws.Range("J58").Value = 0
$ws.Range("K58").Value = 5400
$ws.Range("L58").Value = 0
$ws.Range("M58").Value = -5250
$ws.Range("N58").ClearContents()

$ws.Range("H116").Value = 11530.275
$ws.Range("I116").Value = 12450.048
$ws.Range("J116").Value = 9115.875
$ws.Range("K116").Value = 12450.048
$ws.Range("L116").Value = 9115.875
$ws.Range("M116").Value = -9008.048000000001
$ws.Range("N116").Value = -15999.875

$ws.Range("H138").Value = 2776.6
$ws.Range("J138").Value = 4075.2856
$ws.Range("L138").Value = 12225.8568
$ws.Range("N138").Value = -22505.8568

$ws.Range("H141").Value = 1362.7778
$ws.Range("I141").Value = 1362.7778
$ws.Range("K141").Value = 4088.3334
$ws.Range("M141").Value = 1091.6666

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5751214.5
$ws.Range("I32").Value = 6331755.5
$ws.Range("J32").Value = 18371.5
$ws.Range("K32").Value = 6331755.5
$ws.Range("L32").Value = 18371.5
$ws.Range("M32").Value = -6331468.5
$ws.Range("N32").Value = -18945.5

$ws.Range("H61").Value = 2794272.8
$ws.Range("J61").Value = 4420
$ws.Range("L61").Value = 4420
$ws.Range("N61").Value = -4844

$ws.Range("H97").Value = 856.2857
$ws.Range("I97").Value = 944.6111
$ws.Range("J97").Value = 326.33334
$ws.Range("K97").Value = 944.6111
$ws.Range("L97").Value = 326.33334
$ws.Range("M97").Value = -448.6111
$ws.Range("N97").Value = -1318.33334

$ws.Range("H110").Value = 749.5
$ws.Range("J110").Value = 975
$ws.Range("L110").Value = 975
$ws.Range("N110").Value = -5065

$ws.Range("H122").Value = 3535.0588
$ws.Range("I122").Value = 2093.5
$ws.Range("J122").Value = 4321.364
$ws.Range("K122").Value = 6280.5
$ws.Range("L122").Value = 12964.092
$ws.Range("M122").Value = -3830.5
$ws.Range("N122").Value = -17864.092

$ws.Range("H136").Value = 2794272.8
$ws.Range("J136").Value = 4420
$ws.Range("L136").Value = 13260
$ws.Range("N136").Value = -18360

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2005.7727
$ws.Range("I20").Value = 2175.5715
$ws.Range("J20").Value = 1708.625
$ws.Range("K20").Value = 2175.5715
$ws.Range("L20").Value = 1708.625
$ws.Range("M20").Value = -1928.5715
$ws.Range("N20").Value = -2202.625

$ws.Range("H86").Value = 2041.6666
$ws.Range("I86").Value = 1994.45
$ws.Range("K86").Value = 1994.45
$ws.Range("M86").Value = -871.45

$ws.Range("H89").Value = 2041.6666
$ws.Range("I89").Value = 1994.45
$ws.Range("K89").Value = 9972.25
$ws.Range("M89").Value = -4356.25

$ws.Range("H134").Value = 2764132.5
$ws.Range("I134").Value = 4963469.5
$ws.Range("K134").Value = 14890408.5
$ws.Range("M134").Value = -14887873.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6640.069
$ws.Range("I31").Value = 1899
$ws.Range("K31").Value = 1899
$ws.Range("M31").Value = -1604

$ws.Range("H34").Value = 6640.069
$ws.Range("I34").Value = 1899
$ws.Range("K34").Value = 1899
$ws.Range("M34").Value = -1697

$ws.Range("H58").Value = 3090941.5
$ws.Range("I58").Value = 4116255.2
$ws.Range("K58").Value = 4116255.2
$ws.Range("M58").Value = -4116052.2

$ws.Range("H86").Value = 3800.1667
$ws.Range("J86").Value = 3499.5
$ws.Range("L86").Value = 3499.5
$ws.Range("N86").Value = -5745.5

$ws.Range("H89").Value = 3800.1667
$ws.Range("J89").Value = 3499.5
$ws.Range("L89").Value = 17497.5
$ws.Range("N89").Value = -28729.5

$ws.Range("H136").Value = 3090941.5
$ws.Range("I136").Value = 4116255.2
$ws.Range("K136").Value = 12348765.6
$ws.Range("M136").Value = -12346215.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H107").Value = 869.625
$ws.Range("I107").Value = 581.5714
$ws.Range("K107").Value = 1744.7142
$ws.Range("M107").Value = 175.2857999999999

$ws.Range("H113").Value = 1735.9412
$ws.Range("I113").Value = 1829
$ws.Range("J113").Value = 1697.1666
$ws.Range("K113").Value = 5487
$ws.Range("L113").Value = 5091.4998
$ws.Range("M113").Value = -3317
$ws.Range("N113").Value = -9431.4998

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 183.30435
$ws.Range("J2").Value = 293.81818
$ws.Range("L2").Value = 293.81818
$ws.Range("N2").Value = -519.81818

$ws.Range("H80").Value = 7696.6
$ws.Range("I80").Value = 7085.273
$ws.Range("K80").Value = 7085.273
$ws.Range("M80").Value = -6087.273

$ws.Range("H83").Value = 7696.6
$ws.Range("I83").Value = 7085.273
$ws.Range("K83").Value = 35426.365
$ws.Range("M83").Value = -30434.365

$ws.Range("H102").Value = 2318.16
$ws.Range("I102").Value = 1647.7
$ws.Range("J102").Value = 5000
$ws.Range("K102").Value = 1647.7
$ws.Range("L102").Value = 5000
$ws.Range("M102").Value = -25.70000000000005
$ws.Range("N102").Value = -8244

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2134.926
$ws.Range("I61").Value = 1222.25
$ws.Range("K61").Value = 1222.25
$ws.Range("M61").Value = -1020.25

$ws.Range("H68").Value = 3549.9
$ws.Range("I68").Value = 2500
$ws.Range("K68").Value = 2500
$ws.Range("M68").Value = -1751

$ws.Range("H71").Value = 3549.9
$ws.Range("I71").Value = 2500
$ws.Range("K71").Value = 12500
$ws.Range("M71").Value = -8756

$ws.Range("H100").Value = 7109.409
$ws.Range("I100").Value = 2635.8125
$ws.Range("K100").Value = 2635.8125
$ws.Range("M100").Value = -2094.8125

$ws.Range("H113").Value = 2134.926
$ws.Range("I113").Value = 1222.25
$ws.Range("K113").Value = 1222.25
$ws.Range("M113").Value = 947.75

$ws.Range("H132").Value = 788878.2
$ws.Range("I132").Value = 1050072.5
$ws.Range("K132").Value = 3150217.5
$ws.Range("M132").Value = -3147687.5

$ws.Range("H136").Value = 4083.3547
$ws.Range("I136").Value = 3537.4138
$ws.Range("K136").Value = 10612.2414
$ws.Range("M136").Value = -8062.241399999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 4377958
$ws.Range("I132").Value = 5592427
$ws.Range("K132").Value = 16777281
$ws.Range("M132").Value = -16774751
